$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 14 (old rows 14-16 shift down to 16-18)
$ws.Range("A14:A15").EntireRow.Insert()

# New row 14: Angeleno / Primera (bins, Región Metropolitana)
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44615
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100103
$ws.Range("H14").Value = "Frutos de hueso (carozo)"
$ws.Range("I14").Value = 100103002
$ws.Range("J14").Value = "Ciruela"
$ws.Range("K14").Value = "Angeleno"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 16
$ws.Range("N14").Value = 200000
$ws.Range("O14").Value = 210000
$ws.Range("P14").Value = 205000
$ws.Range("Q14").Value = "`$/bins (450 kilos)"
$ws.Range("R14").Value = "Región Metropolitana"
$ws.Range("S14").Value = 456
$ws.Range("T14").Value = 450

# New row 15: Angeleno / Segunda (bins, Región Metropolitana)
$ws.Range("A15").Value = 2
$ws.Range("B15").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44615
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100103
$ws.Range("H15").Value = "Frutos de hueso (carozo)"
$ws.Range("I15").Value = 100103002
$ws.Range("J15").Value = "Ciruela"
$ws.Range("K15").Value = "Angeleno"
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = 160000
$ws.Range("O15").Value = 170000
$ws.Range("P15").Value = 165000
$ws.Range("Q15").Value = "`$/bins (450 kilos)"
$ws.Range("R15").Value = "Región Metropolitana"
$ws.Range("S15").Value = 367
$ws.Range("T15").Value = 450
